# Refresh the Operation Chaining 7ns Vivado clock-enable power report with the
# latest synthesis run: updated utilization/fanout figures and the removal of
# three stale ap_NS_fsm rows (the FSM signals they referred to no longer appear
# in this run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete ap_NS_fsm[3]/[8]/[4] rows (old rows 14-16); rows below shift up.
$ws.Range("A14:H16").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = [double]"4.179837997071445E-4"

# Row 3
$ws.Range("A3").Value = [double]"1.023563600028865E-4"
$ws.Range("F3").Value = [double]"20.0"

# Row 4
$ws.Range("A4").Value = [double]"1.0136581840924919E-4"
$ws.Range("F4").Value = [double]"28.0"

# Row 5
$ws.Range("A5").Value = [double]"8.802818047115579E-5"
$ws.Range("F5").Value = [double]"16.0"

# Row 6
$ws.Range("A6").Value = [double]"2.899345417972654E-5"
$ws.Range("B6").Value = "firConvolutionOperationChaining_IP/U0/p_pn_reg_118[31]_i_1_n_0"
$ws.Range("C6").Value = [double]"45.09090805053711"
$ws.Range("D6").Value = [double]"15.781818389892578"
$ws.Range("F6").Value = [double]"8.0"
$ws.Range("H6").Value = "FF "

# Row 7
$ws.Range("A7").Value = [double]"2.7669817427522503E-5"
$ws.Range("B7").Value = "firConvolutionOperationChaining_IP/U0/ce0"
$ws.Range("C7").Value = [double]"49.45454406738281"
$ws.Range("D7").Value = [double]"17.30908966064453"
$ws.Range("E7").Value = [double]"18.0"
$ws.Range("F7").Value = [double]"8.0"
$ws.Range("H7").Value = "FF LUT "

# Row 8
$ws.Range("A8").Value = [double]"2.6801455533131957E-5"
$ws.Range("B8").Value = "firConvolutionOperationChaining_IP/U0/shiftRegister_U/firConvolutionOpebkb_ram_U/we0"
$ws.Range("C8").Value = [double]"45.272727966308594"
$ws.Range("D8").Value = [double]"15.80267333984375"
$ws.Range("E8").Value = [double]"32.0"
$ws.Range("F8").Value = [double]"8.0"
$ws.Range("H8").Value = "RAM "

# Row 9
$ws.Range("A9").Value = [double]"2.325745481357444E-5"
$ws.Range("F9").Value = [double]"7.0"

# Row 10
$ws.Range("A10").Value = [double]"6.4512723838561215E-6"
$ws.Range("B10").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[6]"
$ws.Range("F10").Value = [double]"1.0"

# Row 11
$ws.Range("A11").Value = [double]"5.341817995940801E-6"
$ws.Range("B11").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[5]"
$ws.Range("F11").Value = [double]"1.0"

# Row 12
$ws.Range("A12").Value = [double]"4.746545073430752E-6"
$ws.Range("F12").Value = [double]"1.0"

# Row 13
$ws.Range("A13").Value = [double]"2.391272801105515E-6"
$ws.Range("F13").Value = [double]"8.0"

# Row 14
$ws.Range("A14").Value = [double]"5.803636327073036E-7"
$ws.Range("F14").Value = [double]"1.0"
